$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 4658
$ws.Range("F6").Value = 1829
$ws.Range("F7").Value = 44
$ws.Range("F8").Value = 741
$ws.Range("F12").Value = 1132
$ws.Range("F14").Value = 809
$ws.Range("F15").Value = 1262
$ws.Range("F16").Value = 553
$ws.Range("F17").Value = 517
$ws.Range("F19").Value = 174
$ws.Range("F23").Value = 2510
$ws.Range("F25").Value = 1558
$ws.Range("F26").Value = 490
$ws.Range("F29").Value = 4253

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 648
$ws.Range("F4").Value = 14
$ws.Range("F8").Value = 362
$ws.Range("F11").Value = 28
$ws.Range("F24").Value = 51

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 1332
$ws.Range("F7").Value = 256

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1332
$ws.Range("F5").Value = 256
$ws.Range("F8").Value = 14
$ws.Range("F9").Value = 4658
$ws.Range("F12").Value = 1829
$ws.Range("F13").Value = 741
$ws.Range("F14").Value = 362
$ws.Range("F18").Value = 1132
$ws.Range("F20").Value = 28
$ws.Range("F22").Value = 809
$ws.Range("F23").Value = 1262
$ws.Range("F24").Value = 553
$ws.Range("F25").Value = 517
$ws.Range("F27").Value = 174
$ws.Range("F37").Value = 2510
$ws.Range("F39").Value = 51
$ws.Range("F43").Value = 1558
$ws.Range("F44").Value = 490
$ws.Range("F48").Value = 4253
